$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete data rows (old rows 12 and 13) so the table
# shrinks from 12 data rows to 10 data rows (new dimension A1:T11).
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(12).Delete()

# Refresh data rows 2-11 (columns A:T) with the updated TPM figures.

# Row 2
$ws.Range("A2").Value2 = 'ECs'
$ws.Range("B2").Value2 = 'Il34'
$ws.Range("C2").Value2 = 'Ptprz1'
$ws.Range("D2").Value2 = 'FAPs'
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.508158
$ws.Range("H2").Value2 = 1.524474
$ws.Range("I2").Value2 = 0.04862376917091377
$ws.Range("J2").Value2 = 0.05726390893452041
$ws.Range("K2").Value2 = 2
$ws.Range("L2").Value2 = 0.6666666666666666
$ws.Range("M2").Value2 = 0.01638633333333333
$ws.Range("N2").Value2 = 0.049159
$ws.Range("O2").Value2 = 0.02187172081577483
$ws.Range("P2").Value2 = 0.03245268321021395
$ws.Range("Q2").Value2 = 0.008326846374
$ws.Range("R2").Value2 = 0.07494161736600001
$ws.Range("S2").Value2 = 0.001063485504316905
$ws.Range("T2").Value2 = 0.001858367496030531

# Row 3
$ws.Range("A3").Value2 = 'ECs'
$ws.Range("B3").Value2 = 'Il34'
$ws.Range("C3").Value2 = 'Ptprz1'
$ws.Range("D3").Value2 = 'MuSCs'
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.508158
$ws.Range("H3").Value2 = 1.524474
$ws.Range("I3").Value2 = 0.04862376917091377
$ws.Range("J3").Value2 = 0.05726390893452041
$ws.Range("K3").Value2 = 2
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 0.7328155000000001
$ws.Range("N3").Value2 = 1.465631
$ws.Range("O3").Value2 = 0.9781282791842253
$ws.Range("P3").Value2 = 0.9675473167897861
$ws.Range("Q3").Value2 = 0.3723860588490001
$ws.Range("R3").Value2 = 2.234316353094
$ws.Range("S3").Value2 = 0.04756028366659687
$ws.Range("T3").Value2 = 0.05540554143848989

# Row 4
$ws.Range("A4").Value2 = 'FAPs'
$ws.Range("B4").Value2 = 'Il34'
$ws.Range("C4").Value2 = 'Ptprz1'
$ws.Range("D4").Value2 = 'FAPs'
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 5.055996666666666
$ws.Range("H4").Value2 = 15.16799
$ws.Range("I4").Value2 = 0.4837897166804604
$ws.Range("J4").Value2 = 0.5697561244597915
$ws.Range("K4").Value2 = 2
$ws.Range("L4").Value2 = 0.6666666666666666
$ws.Range("M4").Value2 = 0.01638633333333333
$ws.Range("N4").Value2 = 0.049159
$ws.Range("O4").Value2 = 0.02187172081577483
$ws.Range("P4").Value2 = 0.03245268321021395
$ws.Range("Q4").Value2 = 0.08284924671222221
$ws.Range("R4").Value2 = 0.74564322041
$ws.Range("S4").Value2 = 0.01058131361677783
$ws.Range("T4").Value2 = 0.01849011501417285

# Row 5
$ws.Range("A5").Value2 = 'FAPs'
$ws.Range("B5").Value2 = 'Il34'
$ws.Range("C5").Value2 = 'Ptprz1'
$ws.Range("D5").Value2 = 'MuSCs'
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 5.055996666666666
$ws.Range("H5").Value2 = 15.16799
$ws.Range("I5").Value2 = 0.4837897166804604
$ws.Range("J5").Value2 = 0.5697561244597915
$ws.Range("K5").Value2 = 2
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 0.7328155000000001
$ws.Range("N5").Value2 = 1.465631
$ws.Range("O5").Value2 = 0.9781282791842253
$ws.Range("P5").Value2 = 0.9675473167897861
$ws.Range("Q5").Value2 = 3.705112725281666
$ws.Range("R5").Value2 = 22.23067635169
$ws.Range("S5").Value2 = 0.4732084030636827
$ws.Range("T5").Value2 = 0.5512660094456187

# Row 6
$ws.Range("A6").Value2 = 'Inflammatory-Mac'
$ws.Range("B6").Value2 = 'Il34'
$ws.Range("C6").Value2 = 'Ptprz1'
$ws.Range("D6").Value2 = 'FAPs'
$ws.Range("E6").Value2 = 1
$ws.Range("F6").Value2 = 0.3333333333333333
$ws.Range("G6").Value2 = 0.06272999999999999
$ws.Range("H6").Value2 = 0.18819
$ws.Range("I6").Value2 = 0.006002402874876358
$ws.Range("J6").Value2 = 0.007068992335971224
$ws.Range("K6").Value2 = 2
$ws.Range("L6").Value2 = 0.6666666666666666
$ws.Range("M6").Value2 = 0.01638633333333333
$ws.Range("N6").Value2 = 0.049159
$ws.Range("O6").Value2 = 0.02187172081577483
$ws.Range("P6").Value2 = 0.03245268321021395
$ws.Range("Q6").Value2 = 0.00102791469
$ws.Range("R6").Value2 = 0.00925123221
$ws.Range("S6").Value2 = 0.0001312828799030999
$ws.Range("T6").Value2 = 0.0002294077688947045

# Row 7
$ws.Range("A7").Value2 = 'Inflammatory-Mac'
$ws.Range("B7").Value2 = 'Il34'
$ws.Range("C7").Value2 = 'Ptprz1'
$ws.Range("D7").Value2 = 'MuSCs'
$ws.Range("E7").Value2 = 1
$ws.Range("F7").Value2 = 0.3333333333333333
$ws.Range("G7").Value2 = 0.06272999999999999
$ws.Range("H7").Value2 = 0.18819
$ws.Range("I7").Value2 = 0.006002402874876358
$ws.Range("J7").Value2 = 0.007068992335971224
$ws.Range("K7").Value2 = 2
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 0.7328155000000001
$ws.Range("N7").Value2 = 1.465631
$ws.Range("O7").Value2 = 0.9781282791842253
$ws.Range("P7").Value2 = 0.9675473167897861
$ws.Range("Q7").Value2 = 0.045969516315
$ws.Range("R7").Value2 = 0.27581709789
$ws.Range("S7").Value2 = 0.005871119994973259
$ws.Range("T7").Value2 = 0.00683958456707652

# Row 8
$ws.Range("A8").Value2 = 'MuSCs'
$ws.Range("B8").Value2 = 'Il34'
$ws.Range("C8").Value2 = 'Ptprz1'
$ws.Range("D8").Value2 = 'FAPs'
$ws.Range("E8").Value2 = 2
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 4.730544999999999
$ws.Range("H8").Value2 = 9.461089999999999
$ws.Range("I8").Value2 = 0.4526484442488758
$ws.Range("J8").Value2 = 0.355387495084404
$ws.Range("K8").Value2 = 2
$ws.Range("L8").Value2 = 0.6666666666666666
$ws.Range("M8").Value2 = 0.01638633333333333
$ws.Range("N8").Value2 = 0.049159
$ws.Range("O8").Value2 = 0.02187172081577483
$ws.Range("P8").Value2 = 0.03245268321021395
$ws.Range("Q8").Value2 = 0.07751628721833331
$ws.Range("R8").Value2 = 0.4650977233099999
$ws.Range("S8").Value2 = 0.009900200400306229
$ws.Range("T8").Value2 = 0.01153327779484563

# Row 9
$ws.Range("A9").Value2 = 'MuSCs'
$ws.Range("B9").Value2 = 'Il34'
$ws.Range("C9").Value2 = 'Ptprz1'
$ws.Range("D9").Value2 = 'MuSCs'
$ws.Range("E9").Value2 = 2
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 4.730544999999999
$ws.Range("H9").Value2 = 9.461089999999999
$ws.Range("I9").Value2 = 0.4526484442488758
$ws.Range("J9").Value2 = 0.355387495084404
$ws.Range("K9").Value2 = 2
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 0.7328155000000001
$ws.Range("N9").Value2 = 1.465631
$ws.Range("O9").Value2 = 0.9781282791842253
$ws.Range("P9").Value2 = 0.9675473167897861
$ws.Range("Q9").Value2 = 3.4666166994475
$ws.Range("R9").Value2 = 13.86646679779
$ws.Range("S9").Value2 = 0.4427482438485696
$ws.Range("T9").Value2 = 0.3438542172895584

# Row 10
$ws.Range("A10").Value2 = 'Resolving-Mac'
$ws.Range("B10").Value2 = 'Il34'
$ws.Range("C10").Value2 = 'Ptprz1'
$ws.Range("D10").Value2 = 'FAPs'
$ws.Range("E10").Value2 = 2
$ws.Range("F10").Value2 = 0.6666666666666666
$ws.Range("G10").Value2 = 0.093385
$ws.Range("H10").Value2 = 0.280155
$ws.Range("I10").Value2 = 0.008935667024873724
$ws.Range("J10").Value2 = 0.01052347918531281
$ws.Range("K10").Value2 = 2
$ws.Range("L10").Value2 = 0.6666666666666666
$ws.Range("M10").Value2 = 0.01638633333333333
$ws.Range("N10").Value2 = 0.049159
$ws.Range("O10").Value2 = 0.02187172081577483
$ws.Range("P10").Value2 = 0.03245268321021395
$ws.Range("Q10").Value2 = 0.001530237738333333
$ws.Range("R10").Value2 = 0.013772139645
$ws.Range("S10").Value2 = 0.0001954384144707634
$ws.Range("T10").Value2 = 0.0003415151362702371

# Row 11
$ws.Range("A11").Value2 = 'Resolving-Mac'
$ws.Range("B11").Value2 = 'Il34'
$ws.Range("C11").Value2 = 'Ptprz1'
$ws.Range("D11").Value2 = 'MuSCs'
$ws.Range("E11").Value2 = 2
$ws.Range("F11").Value2 = 0.6666666666666666
$ws.Range("G11").Value2 = 0.093385
$ws.Range("H11").Value2 = 0.280155
$ws.Range("I11").Value2 = 0.008935667024873724
$ws.Range("J11").Value2 = 0.01052347918531281
$ws.Range("K11").Value2 = 2
$ws.Range("L11").Value2 = 1
$ws.Range("M11").Value2 = 0.7328155000000001
$ws.Range("N11").Value2 = 1.465631
$ws.Range("O11").Value2 = 0.9781282791842253
$ws.Range("P11").Value2 = 0.9675473167897861
$ws.Range("Q11").Value2 = 0.06843397546750001
$ws.Range("R11").Value2 = 0.410603852805
$ws.Range("S11").Value2 = 0.008740228610402962
$ws.Range("T11").Value2 = 0.01018196404904258
